# Updates cryptos list data (prices, volume %, and some coin name/link/price
# swaps between adjacent rows) to match the latest scrape.
#
# Price (column D) and Volume(1h) (column E) values are stored as plain text
# in the worksheet (e.g. "  +1.25%  ", "66.647.24"), so we force the
# "Text" number format before writing so Excel does not reinterpret
# numeric-looking strings (like "440.67") as floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# Row 2 - Bitcoin
Set-TextValue "D2" "66.682.05"
Set-TextValue "E2" "  +0.93%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.777.15"
Set-TextValue "E3" "  -1.13%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.08%  "

# Row 5 - BNB
Set-TextValue "D5" "440.67"
Set-TextValue "E5" "  +4.36%  "

# Row 6 - Solana
Set-TextValue "D6" "141.93"
Set-TextValue "E6" "  +10.92%  "

# Row 7 - XRP
Set-TextValue "D7" "0.619"
Set-TextValue "E7" "  +2.96%  "

# Row 8 - USDC
Set-TextValue "E8" "  +0.05%  "

# Row 9 - Cardano
Set-TextValue "E9" "  +2.30%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.150"
Set-TextValue "E10" "  -7.89%  "

# Row 11 - ShibaInu
Set-TextValue "E11" "  -10.34%  "

# Row 12 - Avalanche
Set-TextValue "D12" "42.70"
Set-TextValue "E12" "  +6.66%  "

# Row 13 - Polkadot
Set-TextValue "D13" "10.33"
Set-TextValue "E13" "  +4.52%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "4.373.18"
Set-TextValue "E14" "  -1.06%  "

# Row 15 - Uniswap
Set-TextValue "D15" "14.76"
Set-TextValue "E15" "  -7.01%  "

# Row 16 - was WrappedEther, now TRON
Set-TextValue "B16" "TRON"
Set-TextValue "C16" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D16" "0.137"
Set-TextValue "E16" "  -0.36%  "

# Row 17 - was TRON, now WrappedEther
Set-TextValue "B17" "WrappedEther"
Set-TextValue "C17" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D17" "3.774.16"
Set-TextValue "E17" "  -2.31%  "

# Row 18 - Chainlink
Set-TextValue "D18" "19.82"
Set-TextValue "E18" "  +2.29%  "

# Row 19 - Polygon
Set-TextValue "E19" "  +6.77%  "

# Row 20 - WrappedBTC
Set-TextValue "D20" "66.680.71"
Set-TextValue "E20" "  +0.47%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "410.68"
Set-TextValue "E21" "  +2.83%  "

# Row 22 - InternetComputer(DFINITY)
Set-TextValue "D22" "14.45"
Set-TextValue "E22" "  +1.54%  "

# Row 23 - ImmutableX
Set-TextValue "D23" "3.25"
Set-TextValue "E23" "  +9.16%  "

# Row 24 - Litecoin
Set-TextValue "D24" "85.12"
Set-TextValue "E24" "  +1.90%  "

# Row 25 - was PancakeSwap, now EthereumClassic
Set-TextValue "B25" "EthereumClassic"
Set-TextValue "C25" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D25" "36.75"
Set-TextValue "E25" "  -0.17%  "

# Row 26 - was EthereumClassic, now PancakeSwap
Set-TextValue "B26" "PancakeSwap"
Set-TextValue "C26" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D26" "3.38"
Set-TextValue "E26" "  +6.10%  "

# Row 27 - LEO
Set-TextValue "D27" "5.60"
Set-TextValue "E27" "  -2.26%  "

# Row 28 - RenderToken
Set-TextValue "D28" "9.67"
Set-TextValue "E28" "  +33.37%  "

# Row 29 - Filecoin
Set-TextValue "D29" "9.71"
Set-TextValue "E29" "  +4.25%  "

# Row 30 - Bittensor
Set-TextValue "D30" "731.64"
Set-TextValue "E30" "  +6.05%  "

# Row 31 - Cosmos
Set-TextValue "E31" "  +13.39%  "

# Row 32 - Hedera
Set-TextValue "D32" "0.133"
Set-TextValue "E32" "  +11.22%  "

# Row 33 - Toncoin
Set-TextValue "E33" "  -0.11%  "

# Row 34 - InjectiveProtocol
Set-TextValue "D34" "43.15"
Set-TextValue "E34" "  +14.86%  "

# Row 35 - Kaspa
Set-TextValue "E35" "  +4.67%  "

# Row 36 - was NEARProtocol, now OKB
Set-TextValue "B36" "OKB"
Set-TextValue "C36" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D36" "56.21"
Set-TextValue "E36" "  +2.97%  "

# Row 37 - was OKB, now NEARProtocol
Set-TextValue "B37" "NEARProtocol"
Set-TextValue "C37" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D37" "5.54"
Set-TextValue "E37" "  +24.74%  "

# Row 38 - Dai
Set-TextValue "D38" "0.998"
Set-TextValue "E38" "  -0.15%  "

# Row 39 - VeChain
Set-TextValue "D39" "0.0473"
Set-TextValue "E39" "  +5.34%  "

# Row 40 - Fetch.AI
Set-TextValue "D40" "2.69"
Set-TextValue "E40" "  +34.70%  "

# Row 41 - ThetaToken
Set-TextValue "E41" "  -0.97%  "

# Row 42 - was ApeXProtocol, now FirstDigitalUSD
Set-TextValue "B42" "FirstDigitalUSD"
Set-TextValue "C42" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D42" "1.00"
Set-TextValue "E42" "  +0.16%  "

# Row 43 - was FirstDigitalUSD, now ApeXProtocol
Set-TextValue "B43" "ApeXProtocol"
Set-TextValue "C43" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D43" "3.33"
Set-TextValue "E43" "  +7.85%  "

# Row 44 - Stellar
Set-TextValue "E44" "  +4.19%  "

# Row 45 - was PEPE, now TheGraph
Set-TextValue "B45" "TheGraph"
Set-TextValue "C45" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D45" "0.331"
Set-TextValue "E45" "  +16.17%  "

# Row 46 - was TheGraph, now PEPE
Set-TextValue "B46" "PEPE"
Set-TextValue "C46" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D46" "0.0$([char]0x2083)0659"
Set-TextValue "E46" "  -12.93%  "

# Row 47 - LidoDAOToken
Set-TextValue "D47" "3.32"
Set-TextValue "E47" "  +2.71%  "

# Row 48 - was WEMIXToken, now ARBITRUM
Set-TextValue "B48" "ARBITRUM"
Set-TextValue "C48" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D48" "2.08"
Set-TextValue "E48" "  +2.01%  "

# Row 49 - was ARBITRUM, now WEMIXToken
Set-TextValue "B49" "WEMIXToken"
Set-TextValue "C49" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D49" "2.64"
Set-TextValue "E49" "  +4.47%  "

# Row 50 - Monero
Set-TextValue "D50" "142.86"
Set-TextValue "E50" "  -0.68%  "

# Row 51 - Stacks
Set-TextValue "E51" "  +2.90%  "
